$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value2 = "FAPs"
$ws.Cells.Item(2, 2).Value2 = "Lgi3"
$ws.Cells.Item(2, 3).Value2 = "Stx1a"
$ws.Cells.Item(2, 4).Value2 = "ECs"
$ws.Cells.Item(2, 5).Value2 = 3
$ws.Cells.Item(2, 6).Value2 = 1
$ws.Cells.Item(2, 7).Value2 = 1.500469666666667
$ws.Cells.Item(2, 8).Value2 = 4.501409000000001
$ws.Cells.Item(2, 9).Value2 = 0.9943843705197677
$ws.Cells.Item(2, 10).Value2 = 0.9943843705197678
$ws.Cells.Item(2, 11).Value2 = 3
$ws.Cells.Item(2, 12).Value2 = 1
$ws.Cells.Item(2, 13).Value2 = 1.528713333333333
$ws.Cells.Item(2, 14).Value2 = 4.586139999999999
$ws.Cells.Item(2, 15).Value2 = 0.3616028309183848
$ws.Cells.Item(2, 16).Value2 = 0.3616028309183848
$ws.Cells.Item(2, 17).Value2 = 2.293787985695555
$ws.Cells.Item(2, 18).Value2 = 20.64409187126
$ws.Cells.Item(2, 19).Value2 = 0.359572203400944
$ws.Cells.Item(2, 20).Value2 = 0.3595722034009441

# Row 3
$ws.Cells.Item(3, 1).Value2 = "FAPs"
$ws.Cells.Item(3, 2).Value2 = "Lgi3"
$ws.Cells.Item(3, 3).Value2 = "Stx1a"
$ws.Cells.Item(3, 4).Value2 = "FAPs"
$ws.Cells.Item(3, 5).Value2 = 3
$ws.Cells.Item(3, 6).Value2 = 1
$ws.Cells.Item(3, 7).Value2 = 1.500469666666667
$ws.Cells.Item(3, 8).Value2 = 4.501409000000001
$ws.Cells.Item(3, 9).Value2 = 0.9943843705197677
$ws.Cells.Item(3, 10).Value2 = 0.9943843705197678
$ws.Cells.Item(3, 11).Value2 = 3
$ws.Cells.Item(3, 12).Value2 = 1
$ws.Cells.Item(3, 13).Value2 = 1.422612666666667
$ws.Cells.Item(3, 14).Value2 = 4.267837999999999
$ws.Cells.Item(3, 15).Value2 = 0.3365057112737635
$ws.Cells.Item(3, 16).Value2 = 0.3365057112737634
$ws.Cells.Item(3, 17).Value2 = 2.134587153749111
$ws.Cells.Item(3, 18).Value2 = 19.211284383742
$ws.Cells.Item(3, 19).Value2 = 0.334616019881268
$ws.Cells.Item(3, 20).Value2 = 0.334616019881268

# Row 4
$ws.Cells.Item(4, 1).Value2 = "FAPs"
$ws.Cells.Item(4, 2).Value2 = "Lgi3"
$ws.Cells.Item(4, 3).Value2 = "Stx1a"
$ws.Cells.Item(4, 4).Value2 = "M2"
$ws.Cells.Item(4, 5).Value2 = 3
$ws.Cells.Item(4, 6).Value2 = 1
$ws.Cells.Item(4, 7).Value2 = 1.500469666666667
$ws.Cells.Item(4, 8).Value2 = 4.501409000000001
$ws.Cells.Item(4, 9).Value2 = 0.9943843705197677
$ws.Cells.Item(4, 10).Value2 = 0.9943843705197678
$ws.Cells.Item(4, 11).Value2 = 1
$ws.Cells.Item(4, 12).Value2 = 0.3333333333333333
$ws.Cells.Item(4, 13).Value2 = 0.02496166666666666
$ws.Cells.Item(4, 14).Value2 = 0.07488499999999999
$ws.Cells.Item(4, 15).Value2 = 0.005904448619824787
$ws.Cells.Item(4, 16).Value2 = 0.005904448619824787
$ws.Cells.Item(4, 17).Value2 = 0.03745422366277778
$ws.Cells.Item(4, 18).Value2 = 0.337088012965
$ws.Cells.Item(4, 19).Value2 = 0.005871291424090781
$ws.Cells.Item(4, 20).Value2 = 0.005871291424090782

# Row 5
$ws.Cells.Item(5, 1).Value2 = "FAPs"
$ws.Cells.Item(5, 2).Value2 = "Lgi3"
$ws.Cells.Item(5, 3).Value2 = "Stx1a"
$ws.Cells.Item(5, 4).Value2 = "sCs"
$ws.Cells.Item(5, 5).Value2 = 3
$ws.Cells.Item(5, 6).Value2 = 1
$ws.Cells.Item(5, 7).Value2 = 1.500469666666667
$ws.Cells.Item(5, 8).Value2 = 4.501409000000001
$ws.Cells.Item(5, 9).Value2 = 0.9943843705197677
$ws.Cells.Item(5, 10).Value2 = 0.9943843705197678
$ws.Cells.Item(5, 11).Value2 = 3
$ws.Cells.Item(5, 12).Value2 = 1
$ws.Cells.Item(5, 13).Value2 = 1.251315666666667
$ws.Cells.Item(5, 14).Value2 = 3.753947
$ws.Cells.Item(5, 15).Value2 = 0.295987009188027
$ws.Cells.Item(5, 16).Value2 = 0.295987009188027
$ws.Cells.Item(5, 17).Value2 = 1.877561201258111
$ws.Cells.Item(5, 18).Value2 = 16.898050811323
$ws.Cells.Item(5, 19).Value2 = 0.2943248558134649
$ws.Cells.Item(5, 20).Value2 = 0.294324855813465

# Row 6
$ws.Cells.Item(6, 1).Value2 = "sCs"
$ws.Cells.Item(6, 2).Value2 = "Lgi3"
$ws.Cells.Item(6, 3).Value2 = "Stx1a"
$ws.Cells.Item(6, 4).Value2 = "ECs"
$ws.Cells.Item(6, 5).Value2 = 1
$ws.Cells.Item(6, 6).Value2 = 0.3333333333333333
$ws.Cells.Item(6, 7).Value2 = 0.008473666666666666
$ws.Cells.Item(6, 8).Value2 = 0.025421
$ws.Cells.Item(6, 9).Value2 = 0.005615629480232302
$ws.Cells.Item(6, 10).Value2 = 0.005615629480232303
$ws.Cells.Item(6, 11).Value2 = 3
$ws.Cells.Item(6, 12).Value2 = 1
$ws.Cells.Item(6, 13).Value2 = 1.528713333333333
$ws.Cells.Item(6, 14).Value2 = 4.586139999999999
$ws.Cells.Item(6, 15).Value2 = 0.3616028309183848
$ws.Cells.Item(6, 16).Value2 = 0.3616028309183848
$ws.Cells.Item(6, 17).Value2 = 0.01295380721555555
$ws.Cells.Item(6, 18).Value2 = 0.11658426494
$ws.Cells.Item(6, 19).Value2 = 0.002030627517440738
$ws.Cells.Item(6, 20).Value2 = 0.002030627517440739

# Row 7
$ws.Cells.Item(7, 1).Value2 = "sCs"
$ws.Cells.Item(7, 2).Value2 = "Lgi3"
$ws.Cells.Item(7, 3).Value2 = "Stx1a"
$ws.Cells.Item(7, 4).Value2 = "FAPs"
$ws.Cells.Item(7, 5).Value2 = 1
$ws.Cells.Item(7, 6).Value2 = 0.3333333333333333
$ws.Cells.Item(7, 7).Value2 = 0.008473666666666666
$ws.Cells.Item(7, 8).Value2 = 0.025421
$ws.Cells.Item(7, 9).Value2 = 0.005615629480232302
$ws.Cells.Item(7, 10).Value2 = 0.005615629480232303
$ws.Cells.Item(7, 11).Value2 = 3
$ws.Cells.Item(7, 12).Value2 = 1
$ws.Cells.Item(7, 13).Value2 = 1.422612666666667
$ws.Cells.Item(7, 14).Value2 = 4.267837999999999
$ws.Cells.Item(7, 15).Value2 = 0.3365057112737635
$ws.Cells.Item(7, 16).Value2 = 0.3365057112737634
$ws.Cells.Item(7, 17).Value2 = 0.01205474553311111
$ws.Cells.Item(7, 18).Value2 = 0.108492709798
$ws.Cells.Item(7, 19).Value2 = 0.001889691392495486
$ws.Cells.Item(7, 20).Value2 = 0.001889691392495486

# Row 8
$ws.Cells.Item(8, 1).Value2 = "sCs"
$ws.Cells.Item(8, 2).Value2 = "Lgi3"
$ws.Cells.Item(8, 3).Value2 = "Stx1a"
$ws.Cells.Item(8, 4).Value2 = "M2"
$ws.Cells.Item(8, 5).Value2 = 1
$ws.Cells.Item(8, 6).Value2 = 0.3333333333333333
$ws.Cells.Item(8, 7).Value2 = 0.008473666666666666
$ws.Cells.Item(8, 8).Value2 = 0.025421
$ws.Cells.Item(8, 9).Value2 = 0.005615629480232302
$ws.Cells.Item(8, 10).Value2 = 0.005615629480232303
$ws.Cells.Item(8, 11).Value2 = 1
$ws.Cells.Item(8, 12).Value2 = 0.3333333333333333
$ws.Cells.Item(8, 13).Value2 = 0.02496166666666666
$ws.Cells.Item(8, 14).Value2 = 0.07488499999999999
$ws.Cells.Item(8, 15).Value2 = 0.005904448619824787
$ws.Cells.Item(8, 16).Value2 = 0.005904448619824787
$ws.Cells.Item(8, 17).Value2 = 0.0002115168427777777
$ws.Cells.Item(8, 18).Value2 = 0.001903651585
$ws.Cells.Item(8, 19).Value2 = 0.000033157195734005
$ws.Cells.Item(8, 20).Value2 = 0.00003315719573400501

# Row 9
$ws.Cells.Item(9, 1).Value2 = "sCs"
$ws.Cells.Item(9, 2).Value2 = "Lgi3"
$ws.Cells.Item(9, 3).Value2 = "Stx1a"
$ws.Cells.Item(9, 4).Value2 = "sCs"
$ws.Cells.Item(9, 5).Value2 = 1
$ws.Cells.Item(9, 6).Value2 = 0.3333333333333333
$ws.Cells.Item(9, 7).Value2 = 0.008473666666666666
$ws.Cells.Item(9, 8).Value2 = 0.025421
$ws.Cells.Item(9, 9).Value2 = 0.005615629480232302
$ws.Cells.Item(9, 10).Value2 = 0.005615629480232303
$ws.Cells.Item(9, 11).Value2 = 3
$ws.Cells.Item(9, 12).Value2 = 1
$ws.Cells.Item(9, 13).Value2 = 1.251315666666667
$ws.Cells.Item(9, 14).Value2 = 3.753947
$ws.Cells.Item(9, 15).Value2 = 0.295987009188027
$ws.Cells.Item(9, 16).Value2 = 0.295987009188027
$ws.Cells.Item(9, 17).Value2 = 0.01060323185411111
$ws.Cells.Item(9, 18).Value2 = 0.09542908668699999
$ws.Cells.Item(9, 19).Value2 = 0.001662153374562074
$ws.Cells.Item(9, 20).Value2 = 0.001662153374562074

